$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$metaSheet = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$metaSheet.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/claim-item-detail-classification"

# Version: 7.0.0 -> 8.0.0
$metaSheet.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$metaSheet.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$metaSheet.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet 2: Include from Claim Item Detai ---
$includeSheet = $wb.Worksheets.Item("Include from Claim Item Detai")

# System URI: ibm.com -> linuxforhealth.org
$includeSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/claim-item-detail-classification"
